$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Build Diff" worksheet as the last tab, and make it active.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Build Diff"

# ---------------------------------------------------------------------------
# 2. Header row (row 2): labels + sigmoid "Scale" parameter.
#    Formats are copied from existing cells on "Difficulty Scaling" that
#    already carry the desired bold/underline and boxed-and-centred styles,
#    so no new style entries are created in the shared style table.
# ---------------------------------------------------------------------------
$diffScaling = $wb.Worksheets.Item("Difficulty Scaling")

$ws.Range("B2").Value = "Relative Diff"
$ws.Range("C2").Value = "Prob"
$diffScaling.Range("A2").Copy() | Out-Null
$ws.Range("B2:C2").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Value = "Scale"
$ws.Range("F2").Value = 0.9
$diffScaling.Range("L3").Copy() | Out-Null
$ws.Range("E2:F2").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3. Data rows: relative-diff values from -10 to 10, and the sigmoid
#    probability for each. Row 3 uses the scale parameter explicitly; rows
#    4-23 use a plain logistic curve (shared formula).
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = -10
$ws.Range("C3").Formula = "=1/(1+EXP(-`$B3*`$F`$2))"

for ($row = 4; $row -le 23; $row++) {
    $ws.Range("B$row").Value = $row - 13
}
$ws.Range("C4:C23").Formula = "=1/(1+EXP(-B4))"

# ---------------------------------------------------------------------------
# 4. Column widths to fit the new data.
# ---------------------------------------------------------------------------
$ws.Columns("B:C").ColumnWidth = 11.17

# ---------------------------------------------------------------------------
# 5. Selection/view cosmetics on the new sheet.
# ---------------------------------------------------------------------------
$ws.Range("T12").Select() | Out-Null

# ---------------------------------------------------------------------------
# 6. Tidy up stale shared-formula ranges on "Difficulty Scaling" left over
#    from when that sheet had more rows (data now stops at row 53).
# ---------------------------------------------------------------------------
$diffScaling.Range("E4:E53").Formula = $diffScaling.Range("E4").Formula
$diffScaling.Range("B49:B53").Formula = $diffScaling.Range("B49").Formula
$diffScaling.Range("C49:C53").Formula = $diffScaling.Range("C49").Formula

# Re-activate the new sheet so it ends up as the selected/active tab.
$ws.Activate()
